# Updates Jogos_da_Semana_FlashScore_2024-10-31.xlsx:
#  - Refreshes several odds values in rows 2, 4, 5 and 6
#  - Removes the two Switzerland Super League fixtures (rows 7 and 8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (AS Roma - Torino) odds updates
$ws.Range("G2").Value  = 1.75
$ws.Range("I2").Value  = 4.75
$ws.Range("J2").Value  = 2.4
$ws.Range("L2").Value  = 5.5
$ws.Range("U2").Value  = 2.05
$ws.Range("V2").Value  = 1.7
$ws.Range("W2").Value  = 6
$ws.Range("X2").Value  = 7.5
$ws.Range("Z2").Value  = 13
$ws.Range("AA2").Value = 15
$ws.Range("AB2").Value = 34
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 501
$ws.Range("AH2").Value = 10
$ws.Range("AI2").Value = 23
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 9.5
$ws.Range("AS2").Value = 201
$ws.Range("AU2").Value = 9
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 29
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 101
$ws.Range("BB2").Value = 301
$ws.Range("BC2").Value = 151

# Row 4 (Estudiantes L.P. - Ind. Rivadavia) odds update
$ws.Range("G4").Value = 1.57

# Row 5 (Sarmiento Junin - Independiente) odds updates
$ws.Range("H5").Value = 2.88
$ws.Range("I5").Value = 2.2

# Row 6 (Aurora - Independiente) odds update
$ws.Range("G6").Value = 1.53

# Remove the two Switzerland - Super League fixtures entirely
$ws.Rows("7:8").Delete()
